# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on the last data row (row 4)
# of the zh-cn and de-de report sheets to reflect the newly generated handback
# report timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-12 02:30:10"
$zhcn.Range("H4").Value = "2016-03-12 02:30:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-12 02:30:14"
$dede.Range("H4").Value = "2016-03-12 02:30:32"
